$wb = $excel.ActiveWorkbook

# --- GlobalVars: add 3 new global var entries for the new sea creatures ---
$globalVars = $wb.Worksheets.Item("GlobalVars")
$globalVars.Range("A25").Value = "242: Visited the giant whale"
$globalVars.Range("A26").Value = "243: Visited the giant sword fish"
$globalVars.Range("A27").Value = "244: Visited the giant piranha"

# --- Quest - Sea Creatures: update existing coordinate hints and add the new ones ---
$seaCreatures = $wb.Worksheets.Item("Quest - Sea Creatures")
$seaCreatures.Range("A1").Value = "Turtle at 550,402 (map 139)"
$seaCreatures.Range("A2").Value = "Snake at 773,313 (map 112)"
$seaCreatures.Range("A3").Value = "Mermaid at 126,530 (map 163)"
$seaCreatures.Range("A4").Value = "Whale at 173, 773 (map 244)"
$seaCreatures.Range("A5").Value = "Piranha at 311, 446 (map 135)"
$seaCreatures.Range("A6").Value = "Swordfish at 621, 205 (map 77)"
$seaCreatures.Range("A7").Select() | Out-Null

# --- Add new blank worksheet "Tabelle2" after the last sheet (Tabelle3) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Tabelle2"

# --- Todo: add the three new todo entries ---
$todo = $wb.Worksheets.Item("Todo")
$todo.Range("A2").Value = "Finalize Torle's journal text"
$todo.Range("A3").Value = "Finalize cave of the mermaid"
$todo.Range("A4").Value = "Create manyeyes town and castle"
$todo.Range("A5").Select() | Out-Null

# --- Make Todo the active sheet (matches activeTab/tabSelected change in the workbook) ---
$todo.Activate()
